$wb = $excel.ActiveWorkbook

# --- Add sheet "e1" after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$e1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$e1.Name = "e1"

$e1.Range("I8").Value = 12
$e1.Range("J8").Value = 2
$e1.Range("K8").Value = 1
$e1.Range("I9").Value = 23
$e1.Range("J9").Value = 4
$e1.Range("K9").Value = 2
$e1.Range("I10").Value = 34
$e1.Range("J10").Value = 6
$e1.Range("K10").Value = 3
$e1.Range("I11").Value = 45
$e1.Range("J11").Value = 8
$e1.Range("K11").Value = 4
$e1.Range("I12").Value = 56
$e1.Range("J12").Value = 10
$e1.Range("K12").Value = 5
$e1.Range("I13").Value = 67
$e1.Range("J13").Value = 10
$e1.Range("K13").Value = 6

$e1.Columns.Item(9).ColumnWidth = 16.67
$e1.Range("I8:K13").Select()

# --- Add sheet "e1_shifted" after "e1" ---
$e1s = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $e1)
$e1s.Name = "e1_shifted"

$e1s.Range("C7").Value = 12
$e1s.Range("D7").Value = 2
$e1s.Range("E7").Value = 1
$e1s.Range("C8").Value = 23
$e1s.Range("D8").Value = 4
$e1s.Range("E8").Value = 2
$e1s.Range("C9").Value = 34
$e1s.Range("D9").Value = 6
$e1s.Range("E9").Value = 3
$e1s.Range("C10").Value = 45
$e1s.Range("D10").Value = 8
$e1s.Range("E10").Value = 4
$e1s.Range("C11").Value = 56
$e1s.Range("D11").Value = 10
$e1s.Range("E11").Value = 5
$e1s.Range("C12").Value = 67
$e1s.Range("D12").Value = 10
$e1s.Range("E12").Value = 6

$e1s.Columns.Item(3).ColumnWidth = 16.67
$e1s.Range("C7").Select()

# --- Add sheet "e2" after "e1_shifted" ---
$e2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $e1s)
$e2.Name = "e2"

$e2.Range("B3").Value = "Heading 1"
$e2.Range("C3").Value = "Heading 2"
$e2.Range("D3").Value = "Heading 3"
$e2.Range("B4").Value = 12
$e2.Range("C4").Value = 2
$e2.Range("D4").Value = 1
$e2.Range("B5").Value = 23
$e2.Range("C5").Value = 4
$e2.Range("D5").Value = 2
$e2.Range("B6").Value = 34
$e2.Range("C6").Value = 6
$e2.Range("D6").Value = 3
$e2.Range("B7").Value = 45
$e2.Range("C7").Value = 8
$e2.Range("D7").Value = 4
$e2.Range("B8").Value = 56
$e2.Range("C8").Value = 10
$e2.Range("D8").Value = 5
$e2.Range("B9").Value = 67
$e2.Range("C9").Value = 10
$e2.Range("D9").Value = 6

$e2.Range("B3:D9").Select()

# --- Add sheet "e2_shifted" after "e2" ---
$e2s = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $e2)
$e2s.Name = "e2_shifted"

$e2s.Range("B4").Value = "Heading 1"
$e2s.Range("C4").Value = "Heading 2"
$e2s.Range("D4").Value = "Heading 3"
$e2s.Range("B5").Value = 12
$e2s.Range("C5").Value = 2
$e2s.Range("D5").Value = 1
$e2s.Range("B6").Value = 23
$e2s.Range("C6").Value = 4
$e2s.Range("D6").Value = 2
$e2s.Range("B7").Value = 34
$e2s.Range("C7").Value = 6
$e2s.Range("D7").Value = 3
$e2s.Range("B8").Value = 45
$e2s.Range("C8").Value = 8
$e2s.Range("D8").Value = 4
$e2s.Range("B9").Value = 56
$e2s.Range("C9").Value = 10
$e2s.Range("D9").Value = 5
$e2s.Range("B10").Value = 67
$e2s.Range("C10").Value = 10
$e2s.Range("D10").Value = 6

$e2s.Range("H15").Select()

# Activate e2 as the active tab (index 7, 0-based) at the end
$e2.Activate()

Write-Host "done"
